$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update percentages on existing parties (rows 2-8) ---
$ws1.Range("C2").Value = 0.225
$ws1.Range("C3").Value = 0.218
$ws1.Range("C4").Value = 0.145
$ws1.Range("C5").Value = 0.083
$ws1.Range("C6").Value = 0.107
$ws1.Range("C7").Value = 0.045
$ws1.Range("C8").Value = 0.026

# --- Row 10 (Europa Verde - Sinistra Italiana): abbreviation + percentage update ---
$ws1.Range("B10").Value = "EV - SI"
$ws1.Range("C10").Value = 0.043

# --- Row 11 becomes a real party/list row ("Art. 1 - MDP") replacing "Altro 1" ---
$ws1.Range("A11").Value = "Art. 1 - MDP"
$ws1.Range("B11").Value = "MDP"
$ws1.Range("C11").Value = 0.019
$ws1.Range("D11").Value = "SX"
$ws1.Range("I11").Value = 15
$ws1.Range("J11").Value = "CENTRO"

# --- Remove the old "Altro 2".."Altro 5" rows (12-15) entirely ---
$ws1.Rows.Item(12).Resize(4, 1).EntireRow.Delete() | Out-Null

# --- Selection on sheet1 moves to C11 ---
$ws1.Range("C11").Select() | Out-Null

# --- Add a new sheet "altri_dati" after "liste_naz" ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "altri_dati"
$ws2.Range("A1").Value = "Astensione"
$ws2.Range("A2").Value = 0.5
$ws2.Columns.Item(1).NumberFormat = "0%"
$ws2.Columns.Item(1).ColumnWidth = 11
$ws2.Range("B1").Select() | Out-Null

# --- Re-activate the first sheet so it stays the selected tab ---
$ws1.Activate()
